# ConstraintGenerationDevelopperDoc.docx edit script
# Bug 382733: [validation] Papyrus shall to compute constraints from the
# profile at the model level

$d = $word.ActiveDocument

function Find-ParagraphByExactText($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Insert the two new paragraphs right before the "Constraint in OCL"
#    Heading 2 paragraph (after "Generate constraints as EMF validation
#    plugins").
# ---------------------------------------------------------------------
$target = Find-ParagraphByExactText("Constraint in OCL`r")
$r = $target.Range
$r.Collapse(1)
$r.InsertBefore("ZZZ_PLACEHOLDER_1`r")

$ph = Find-ParagraphByExactText("ZZZ_PLACEHOLDER_1`r")
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Constraint written in </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>OCL</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> in the profile can be generated into the definition of the profile and taken in charge during the validation of the model.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The problem w</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>it</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>h</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> this use case </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">is </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>to</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>know if it is possible to respect the requirement ValidationReq002</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$ph.Range.InsertXML($xml1)

Write-Output "step1 done"

# ---------------------------------------------------------------------
# 2) Insert the two new paragraphs right before the "Constraint in Java"
#    Heading 2 paragraph (after the blank paragraph that follows
#    "Constraint in OCL").
# ---------------------------------------------------------------------
$target2 = Find-ParagraphByExactText("Constraint in Java`r")
$r2 = $target2.Range
$r2.Collapse(1)
$r2.InsertBefore("ZZZ_PLACEHOLDER_2`r")

$ph2 = Find-ParagraphByExactText("ZZZ_PLACEHOLDER_2`r")
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The user can generate plugins that wrap constraint and can be used in the EMF plugin validation.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The constraint can be generated in Java code, or directly from OCL</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r><w:bookmarkStart w:id="100" w:name="_GoBack"/><w:bookmarkEnd w:id="100"/></w:p>'
$ph2.Range.InsertXML($xml2)

Write-Output "step2 done"
